$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Category_Num for the "Passive" row (row 5) from 0 to 4
$ws.Range("D5").Value = 4

# Match the resulting active cell / selection recorded in the saved file
$ws.Range("D5").Select()
